# Remove the trailing "footer" paragraphs that were scraped along with the
# course description: the blank paragraph, the "Ver no Jupiter ..." line,
# and the "© 2020 ..." copyright line that follow the Requisitos entry
# "LOM3015: Termodinâmica de Materiais (Requisito fraco)".
$d = $word.ActiveDocument

# Locate the paragraph holding the "Ver no Jupiter ..." text; the blank
# paragraph right before it and the copyright paragraph right after it are
# removed together with it.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $target = $p
        break
    }
}

$prev = $target.Previous(1)
$next = $target.Next(1)

$startRange = $prev.Range.Start
$endRange = $next.Range.End

$r = $d.Range($startRange, $endRange)
$r.Delete()
